$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$r = $ws.Range("Z5")
$r.Value = 99
$r.BorderAround(1)
